$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-15 Thursday", "2024-08-16 Friday"),
    @("848×6=", "973×8="),
    @("157×3=", "744×5="),
    @("155×3=", "423×4="),
    @("204×6=", "802×4="),
    @("611×4=", "463×3="),
    @("682×3=", "110×5="),
    @("483×3=", "494×7="),
    @("478×4=", "520×7="),
    @("234×9=", "449×5="),
    @("746×3=", "631×4="),
    @("443×3=", "200×7="),
    @("106×2=", "953×9="),
    @("548×6=", "658×7="),
    @("301×6=", "768×6="),
    @("424×2=", "293×5="),
    @("445×8=", "962×4="),
    @("370×2=", "462×8="),
    @("189×5=", "267×9="),
    @("906×7=", "332×7="),
    @("762×8=", "768×3="),
    @("387×7=", "657×9="),
    @("843×3=", "489×6="),
    @("795×4=", "940×8="),
    @("158×7=", "937×3="),
    @("665×3=", "777×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
